$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: cell address, new value (always written as text to preserve
# formatting such as leading/trailing zeros, exactly like the source data).
$changes = @(
    @("D2", "271.08"),
    @("D3", "22.79"),
    @("D4", "6.345"),
    @("D5", "0.06202"),
    @("D7", "6.696"),
    @("D8", "1.387"),
    @("D9", "0.8309"),
    @("D11", "0.1604"),
    @("D12", "0.08290"),
    @("D14", "0.03181"),
    @("B15", "BitMartToken"),
    @("C15", "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"),
    @("D15", "0.09333"),
    @("E15", "14BitMartTokenBMX"),
    @("B16", "MCDex"),
    @("C16", "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"),
    @("D16", "3.855"),
    @("E16", "15MCDexMCB"),
    @("B17", "BitForexToken"),
    @("C17", "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"),
    @("D17", "0.001656"),
    @("E17", "16BitForexTokenBF"),
    @("B18", "CoinExToken"),
    @("C18", "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"),
    @("D18", "0.04732"),
    @("E18", "17CoinExTokenCET"),
    @("B19", "TigerCash"),
    @("C19", "https://coinranking.com/coin/6hIn06L2+tigercash-tch"),
    @("D19", "0.006303"),
    @("E19", "18TigerCashTCH"),
    @("B20", "HotbitToken"),
    @("C20", "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"),
    @("D20", "0.005663"),
    @("E20", "19HotbitTokenHTB"),
    @("B21", "BitKan"),
    @("C21", "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"),
    @("D21", "0.001078"),
    @("E21", "20BitKanKAN"),
    @("B22", "NitroEx"),
    @("C22", "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"),
    @("D22", "0.0001501"),
    @("E22", "21NitroExNTX"),
    @("B23", "LEO"),
    @("C23", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"),
    @("D23", "3.718"),
    @("E23", "22LEOLEO"),
    @("B24", "BTSEToken"),
    @("C24", "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"),
    @("D24", "2.399"),
    @("E24", "23BTSETokenBTSE"),
    @("B25", "BitpandaEcosystemToken"),
    @("C25", "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"),
    @("D25", "0.3348"),
    @("E25", "24BitpandaEcosystemTokenBEST"),
    @("B26", "ProBitToken"),
    @("C26", "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"),
    @("D26", "0.1240"),
    @("E26", "25ProBitTokenPROB"),
    @("D27", "0.0002706"),
    @("D40", "0.04684"),
    @("D41", "0.007019"),
    @("D42", "0.1160"),
    @("D43", "0.003293"),
    @("E43", "42CEJICEJIWorstin24h"),
    @("D44", "0.01167"),
    @("D45", "0.00006264"),
    @("D46", "0.0009906"),
    @("D48", "0.9206"),
    @("D49", "0.002120"),
    @("E50", "49CryptobidCoinCBC")
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $newValue = $change[1]
    $rng = $ws.Range($cellRef)
    # Preserve the existing style while forcing a Text number format so Excel
    # does not reinterpret numeric-looking strings (e.g. "271.08") as numbers.
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}
